$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 28: add the new activity description in column A and hours (1.5) in column B.
$ws.Range("A28").Value = "3. iterace - sekvenční diagramy a jeden scénář"
$ws.Range("B28").Value = 1.5

# Move the active selection to A29 (was B19).
$ws.Range("A29").Select()
